$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '43.814.76'
$ws.Cells.Item(2, 5).Value = '  -0.38%  '
$ws.Cells.Item(3, 4).Value = '2.345.83'
$ws.Cells.Item(3, 5).Value = '  -0.33%  '
$ws.Cells.Item(4, 5).Value = '  +0.00%  '
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '239.05'
$cell.ClearFormats()
$ws.Cells.Item(5, 5).Value = '  -1.08%  '
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.666'
$cell.ClearFormats()
$ws.Cells.Item(6, 5).Value = '  -4.14%  '
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = '72.72'
$cell.ClearFormats()
$ws.Cells.Item(7, 5).Value = '  -4.68%  '
$ws.Cells.Item(8, 5).Value = '  -0.02%  '
$ws.Cells.Item(9, 5).Value = '  -5.64%  '
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.101'
$cell.ClearFormats()
$ws.Cells.Item(10, 5).Value = '  -1.19%  '
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '60.11'
$cell.ClearFormats()
$ws.Cells.Item(11, 5).Value = '  +4.76%  '
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = '32.72'
$cell.ClearFormats()
$ws.Cells.Item(12, 5).Value = '  -1.82%  '
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.108'
$cell.ClearFormats()
$ws.Cells.Item(13, 5).Value = '  +0.02%  '
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.26'
$cell.ClearFormats()
$ws.Cells.Item(14, 5).Value = '  -3.24%  '
$ws.Cells.Item(15, 4).Value = '2.693.82'
$ws.Cells.Item(15, 5).Value = '  -0.46%  '
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '16.07'
$cell.ClearFormats()
$ws.Cells.Item(16, 5).Value = '  -4.40%  '
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.897'
$cell.ClearFormats()
$ws.Cells.Item(17, 5).Value = '  -3.30%  '
$ws.Cells.Item(18, 4).Value = '2.341.74'
$ws.Cells.Item(18, 5).Value = '  -0.60%  '
$ws.Cells.Item(19, 4).Value = '43.732.07'
$ws.Cells.Item(19, 5).Value = '  -0.43%  '
$ws.Cells.Item(20, 5).Value = '  +0.04%  '
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.67'
$cell.ClearFormats()
$ws.Cells.Item(21, 5).Value = '  -0.08%  '
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '78.48'
$cell.ClearFormats()
$ws.Cells.Item(22, 5).Value = '  +0.94%  '
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '251.50'
$cell.ClearFormats()
$ws.Cells.Item(23, 5).Value = '  -4.31%  '
$ws.Cells.Item(24, 5).Value = '  +0.19%  '
$ws.Cells.Item(25, 5).Value = '  +3.07%  '
$ws.Cells.Item(26, 5).Value = '  +1.50%  '
$ws.Cells.Item(27, 5).Value = '  -1.74%  '
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = '10.39'
$cell.ClearFormats()
$ws.Cells.Item(28, 5).Value = '  -4.98%  '
$ws.Cells.Item(29, 5).Value = '  -1.96%  '
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = '176.65'
$cell.ClearFormats()
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = '22.22'
$cell.ClearFormats()
$ws.Cells.Item(31, 5).Value = '  -3.94%  '
$ws.Cells.Item(32, 5).Value = '  -1.16%  '
$ws.Cells.Item(33, 5).Value = '  -2.96%  '
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0744'
$cell.ClearFormats()
$ws.Cells.Item(34, 5).Value = '  -2.38%  '
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.07'
$cell.ClearFormats()
$ws.Cells.Item(35, 5).Value = '  -5.98%  '
$ws.Cells.Item(36, 5).Value = '  -2.30%  '
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.74'
$cell.ClearFormats()
$ws.Cells.Item(37, 5).Value = '  -2.20%  '
$ws.Cells.Item(38, 5).Value = '  -0.48%  '
$ws.Cells.Item(39, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.37'
$cell.ClearFormats()
$ws.Cells.Item(39, 5).Value = '  -2.25%  '
$ws.Cells.Item(40, 2).Value = 'FTXToken'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.79'
$cell.ClearFormats()
$ws.Cells.Item(40, 5).Value = '  +23.53%  '
$ws.Cells.Item(41, 5).Value = '  -4.30%  '
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = '65.18'
$cell.ClearFormats()
$ws.Cells.Item(42, 5).Value = '  +15.22%  '
$ws.Cells.Item(43, 5).Value = '  +0.71%  '
$ws.Cells.Item(44, 5).Value = '  -2.08%  '
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '18.74'
$cell.ClearFormats()
$ws.Cells.Item(45, 5).Value = '  -2.78%  '
$ws.Cells.Item(46, 5).Value = '  -11.73%  '
$ws.Cells.Item(47, 5).Value = '  +0.00%  '
$ws.Cells.Item(48, 5).Value = '  -2.91%  '
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.15'
$cell.ClearFormats()
$ws.Cells.Item(49, 5).Value = '  -3.23%  '
$ws.Cells.Item(50, 2).Value = 'Aave'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '97.97'
$cell.ClearFormats()
$ws.Cells.Item(50, 5).Value = '  -3.87%  '
$ws.Cells.Item(51, 2).Value = 'NEARProtocol'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.40'
$cell.ClearFormats()
$ws.Cells.Item(51, 5).Value = '  -5.21%  '
